# Generate Report for Archive
# The "ebbbe0f4-60d8-47f2-97a4-130aaff75ada" entry has moved up in the report
# (its status became "In Translation") and now sits at row 5, pushing the
# "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8" and "4d93103c-a919-4d01-b99a-dbe0d34ebbde"
# rows down to rows 6 and 7 respectively, on every sheet of the workbook.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$cellAddr,
        [string]$newValue,
        [bool]$isHyperlink
    )

    $targetRange = $ws.Range($cellAddr)
    $targetAddress = $targetRange.Address()
    $targetRange.Value2 = $newValue

    if ($isHyperlink) {
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $targetAddress) {
                $hl.TextToDisplay = $newValue
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet: columns A (File Name, hyperlinked), B (zh-cn status),
# C (de-de status)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $wsOverview "A5" "ebbbe0f4-60d8-47f2-97a4-130aaff75ada.md" $true
$wsOverview.Range("B5").Value2 = "In Translation"
$wsOverview.Range("C5").Value2 = "In Translation"

Set-CellAndHyperlink $wsOverview "A6" "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md" $true
$wsOverview.Range("B6").Value2 = "In Translation"
$wsOverview.Range("C6").Value2 = "In Translation"

Set-CellAndHyperlink $wsOverview "A7" "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md" $true
$wsOverview.Range("B7").Value2 = "Ready for handoff"
$wsOverview.Range("C7").Value2 = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: A=Source File Name (hyperlinked), B=Status,
# C=Latest Handoff File (hyperlinked), D=Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A5" "ebbbe0f4-60d8-47f2-97a4-130aaff75ada.md" $true
$wsZhCn.Range("B5").Value2 = "In Translation"
Set-CellAndHyperlink $wsZhCn "C5" "ebbbe0f4-60d8-47f2-97a4-130aaff75ada.8ec8dff7ea39a4a8660d34cea34430c024d49e12.zh-cn.xlf" $true
$wsZhCn.Range("D5").Value2 = "2016-03-02 09:44:20"

Set-CellAndHyperlink $wsZhCn "A6" "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md" $true
$wsZhCn.Range("B6").Value2 = "In Translation"
Set-CellAndHyperlink $wsZhCn "C6" "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.zh-cn.xlf" $true
$wsZhCn.Range("D6").Value2 = "2016-03-02 09:38:32"

Set-CellAndHyperlink $wsZhCn "A7" "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md" $true
$wsZhCn.Range("B7").Value2 = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "C7" "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.zh-cn.xlf" $true
$wsZhCn.Range("D7").Value2 = "2016-03-02 09:39:16"

# ---------------------------------------------------------------------------
# de-de sheet: same layout as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A5" "ebbbe0f4-60d8-47f2-97a4-130aaff75ada.md" $true
$wsDeDe.Range("B5").Value2 = "In Translation"
Set-CellAndHyperlink $wsDeDe "C5" "ebbbe0f4-60d8-47f2-97a4-130aaff75ada.8ec8dff7ea39a4a8660d34cea34430c024d49e12.de-de.xlf" $true
$wsDeDe.Range("D5").Value2 = "2016-03-02 09:44:31"

Set-CellAndHyperlink $wsDeDe "A6" "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md" $true
$wsDeDe.Range("B6").Value2 = "In Translation"
Set-CellAndHyperlink $wsDeDe "C6" "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.de-de.xlf" $true
$wsDeDe.Range("D6").Value2 = "2016-03-02 09:38:42"

Set-CellAndHyperlink $wsDeDe "A7" "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md" $true
$wsDeDe.Range("B7").Value2 = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "C7" "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.de-de.xlf" $true
$wsDeDe.Range("D7").Value2 = "2016-03-02 09:39:27"

$wb.Save()
